# Insert a new pair of rows (a new "Primera"/"Segunda" weekly price record)
# right before row 446, pushing all the existing records down by two rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 446 (existing row 446 and everything
# below shifts down to 448+).
$ws.Rows.Item(446).Resize(2).Insert()

# New row 446: "Primera" quality entry for the new weekly date (2023-09-11 -> serial 45180).
$ws.Range("A446").Value = 1
$ws.Range("B446").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C446").Value = "Arica y Parinacota"
$ws.Range("D446").Value = 45180
$ws.Range("E446").Value = 15
$ws.Range("F446").Value = 100112043
$ws.Range("G446").Value = "Pepino ensalada"
$ws.Range("H446").Value = "Sin especificar"
$ws.Range("I446").Value = "Primera"
$ws.Range("J446").Value = 150
$ws.Range("K446").Value = 11000
$ws.Range("L446").Value = 12000
$ws.Range("M446").Value = 11500
$ws.Range("N446").Value = "$/caja 70 unidades"
$ws.Range("O446").Value = "Región de Arica y Parinacota"
$ws.Range("P446").Value = 164
$ws.Range("Q446").Value = 70
$ws.Range("R446").Value = "Hortaliza"

# New row 447: "Segunda" quality entry for the same new weekly date.
$ws.Range("A447").Value = 1
$ws.Range("B447").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C447").Value = "Arica y Parinacota"
$ws.Range("D447").Value = 45180
$ws.Range("E447").Value = 15
$ws.Range("F447").Value = 100112043
$ws.Range("G447").Value = "Pepino ensalada"
$ws.Range("H447").Value = "Sin especificar"
$ws.Range("I447").Value = "Segunda"
$ws.Range("J447").Value = 150
$ws.Range("K447").Value = 9000
$ws.Range("L447").Value = 10000
$ws.Range("M447").Value = 9500
$ws.Range("N447").Value = "$/caja 100 unidades"
$ws.Range("O447").Value = "Región de Arica y Parinacota"
$ws.Range("P447").Value = 95
$ws.Range("Q447").Value = 100
$ws.Range("R447").Value = "Hortaliza"

# Make sure the date columns keep the workbook's datetime display format.
$ws.Range("D446:D447").NumberFormat = "YYYY-MM-DD HH:MM:SS"
